# Adds a new "2021" column (column R) to the statistics table on the
# active sheet, mirroring the existing "2020" column (Q) for layout/
# formatting, then overwrites the three data cells with the 2021 figures.
#
#   R2 -> blank divider cell (same style as Q2)
#   R3 -> 2021                (year header, same style as Q3)
#   R4 -> 202551               (population count, same style as Q4)
#   R5 -> 2.9794303052841493   (percentage, same style as Q5)
#
# Finally the active selection is moved to R2, matching the saved
# workbook's view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column Q's formatting/layout into the new column R for every
# populated row so the new column visually matches the rest of the table
# (borders, fonts, number formats, row-bottom rule, etc.).
$ws.Range("Q2").Copy($ws.Range("R2"))
$ws.Range("Q3").Copy($ws.Range("R3"))
$ws.Range("Q4").Copy($ws.Range("R4"))
$ws.Range("Q5").Copy($ws.Range("R5"))

# Overwrite the copied values with the real 2021 figures.
$ws.Range("R3").Value = 2021
$ws.Range("R4").Value = 202551
$ws.Range("R5").Value = 2.9794303052841493

# Match the workbook's saved selection state.
$ws.Range("R2").Select()
